$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-10 Wednesday" "2024-07-11 Thursday"

Replace-Text "55×14=" "42×17="
Replace-Text "40×16=" "38×90="
Replace-Text "32×65=" "22×30="
Replace-Text "89×29=" "35×73="
Replace-Text "52×99=" "95×46="
Replace-Text "31×22=" "19×57="
Replace-Text "67×26=" "73×16="
Replace-Text "97×71=" "24×88="
Replace-Text "46×16=" "14×67="
Replace-Text "99×38=" "84×82="
Replace-Text "79×16=" "25×87="
Replace-Text "35×30=" "49×64="
Replace-Text "72×32=" "53×80="
Replace-Text "59×95=" "83×73="
Replace-Text "77×16=" "54×98="
Replace-Text "58×49=" "25×35="
Replace-Text "99×23=" "52×32="
Replace-Text "67×40=" "91×37="
Replace-Text "92×92=" "12×17="
Replace-Text "28×71=" "87×31="
Replace-Text "84×21=" "66×48="
Replace-Text "49×34=" "96×19="
Replace-Text "66×30=" "54×17="
Replace-Text "72×75=" "71×73="
Replace-Text "57×54=" "11×28="
